$d = $word.ActiveDocument

# Locate the table that holds the inspection fields (the one with the
# "Medidor" / "${medidor}" row) and find that row's index.
$targetTable = $null
$medidorRowIndex = -1
foreach ($tbl in $d.Tables) {
    for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
        $labelCell = $tbl.Cell($i, 1)
        $labelText = $labelCell.Range.Text.TrimEnd([char]13, [char]7)
        if ($labelText -eq "Medidor") {
            $targetTable = $tbl
            $medidorRowIndex = $i
            break
        }
    }
    if ($targetTable -ne $null) { break }
}

if ($targetTable -eq $null) {
    throw "Could not find the 'Medidor' row"
}

# Insert a brand-new row right below the "Medidor" row by inserting it
# above the row that currently follows "Medidor" (Word's Rows.Add(before)
# inserts immediately above the reference row).
$referenceRow = $targetTable.Rows.Item($medidorRowIndex + 1)
$newRow = $targetTable.Rows.Add($referenceRow)

$newRow.Cells.Item(1).Range.Text = "Medidor Anomalia"
$newRow.Cells.Item(2).Range.Text = "`${medidor_anomalia}"
